$wb = $excel.ActiveWorkbook

# --- Update "Logs" sheet: append new row 10 ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A10").Value = "Demo inplannen"
$logs.Range("B10").Value = "klantenservice@testbedrijf123.nl"
$logs.Range("C10").Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$logs.Range("D10").Value = "Intern verzoek / Actie voor medewerker"
$logs.Range("E10").Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$logs.Range("F10").Value = "2025-08-13 21:09:38"
$logs.Range("G10").Value = "Nee"
$logs.Range("H10").Value = "Ja"
$logs.Range("I10").Value = "Nee"
$logs.Range("J10").Value = "Nee"

# --- Extend conditional formatting ranges to include the new row ---
$logs.Range("D2:D9").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D10"))
$logs.Range("G2:G9").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G10"))
$logs.Range("H2:H9").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H10"))
$logs.Range("I2:I9").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I10"))
$logs.Range("J2:J9").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J10"))

# --- Update "Dashboard" sheet: bump the count for this category ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 9
